$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 into I1:J1, then set header text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row, I0, IF
$data = @(
    @(2,7,7),
    @(3,7,7),
    @(4,8,8),
    @(5,6,6),
    @(6,7,8),
    @(7,9,9),
    @(8,8,8),
    @(9,9,9),
    @(10,7,7),
    @(11,9,9),
    @(12,8,8),
    @(13,9,9),
    @(14,8,9),
    @(15,9,9),
    @(16,8,8),
    @(17,8,8),
    @(18,9,9),
    @(19,8,8),
    @(20,9,9),
    @(21,7,7),
    @(22,10,10),
    @(23,9,9),
    @(24,7,7),
    @(25,8,8),
    @(26,8,8),
    @(27,7,7),
    @(28,6,6),
    @(29,7,7),
    @(30,7,7),
    @(31,8,8),
    @(32,7,8),
    @(33,7,7),
    @(34,6,7),
    @(35,7,7),
    @(36,7,7),
    @(37,7,7),
    @(38,6,6),
    @(39,6,6),
    @(40,6,7),
    @(41,10,10),
    @(42,1,1),
    @(43,7,7),
    @(44,1,2),
    @(45,8,8),
    @(46,6,6),
    @(47,6,7),
    @(48,9,9),
    @(49,1,2),
    @(50,1,2),
    @(51,6,6),
    @(52,9,9),
    @(53,7,7),
    @(54,6,7),
    @(55,9,9),
    @(56,7,7),
    @(57,6,7),
    @(58,4,6),
    @(59,5,6),
    @(60,6,6),
    @(61,6,6),
    @(62,6,6),
    @(63,8,8),
    @(64,9,9),
    @(65,8,8),
    @(66,7,7),
    @(67,6,7),
    @(68,9,9),
    @(69,4,5),
    @(70,6,6),
    @(71,5,6),
    @(72,5,6)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 9).Value = $item[1]
    $ws.Cells.Item($r, 10).Value = $item[2]
}
